$wb = $excel.ActiveWorkbook

# --- Rename "Story Map Data" -> "Feature Data" ---
$wsFeature = $wb.Worksheets.Item("Story Map Data")
$wsFeature.Name = "Feature Data"

# --- Populate the "Feature Data" sheet ---
# Header row: FEATURE_ID | TITLE | US_ID
# Values are written in an order chosen so that the shared-string table
# ends up built in the same sequence as the target workbook.
$wsFeature.Range("C1").Value = "US_ID"

$wsFeature.Range("A2").Value = "Feature 1"
$wsFeature.Range("B2").Value = "Feature 1 TITLE"
$wsFeature.Range("A10").Value = "Feature 2"
$wsFeature.Range("B10").Value = "Feature 2 TITLE"
$wsFeature.Range("A16").Value = "Feature 3"
$wsFeature.Range("B16").Value = "Feature 3 TITLE"

$wsFeature.Range("A1").Value = "FEATURE_ID"
$wsFeature.Range("B1").Value = "TITLE"

for ($r = 2; $r -le 9; $r++) {
    $wsFeature.Range("A$r").Value = "Feature 1"
    $wsFeature.Range("B$r").Value = "Feature 1 TITLE"
    $wsFeature.Range("C$r").Value = $r - 1
}
for ($r = 10; $r -le 15; $r++) {
    $wsFeature.Range("A$r").Value = "Feature 2"
    $wsFeature.Range("B$r").Value = "Feature 2 TITLE"
    $wsFeature.Range("C$r").Value = $r - 1
}
for ($r = 16; $r -le 19; $r++) {
    $wsFeature.Range("A$r").Value = "Feature 3"
    $wsFeature.Range("B$r").Value = "Feature 3 TITLE"
    $wsFeature.Range("C$r").Value = $r - 1
}

# --- Page setup for the "Feature Data" sheet ---
$wsFeature.PageSetup.PaperSize = 9
$wsFeature.PageSetup.Orientation = 1

# --- Selection / active sheet ---
# Selecting A2 on "Feature Data" makes it the active sheet & cell, which
# also clears the previous selection state (tabSelected/topLeftCell) on
# the "US Data" sheet.
$wsFeature.Range("A2").Select() | Out-Null

Write-Host "done"
